# feat: add missing strings
#
# 1. Update the Chinese homepage meta description (B3) to the new, shorter copy.
# 2. Add the missing English homepage meta description (C3), matching the
#    existing style used by sibling data cells in that row (copy format from C4).
# 3. Append a brand-new "aboutPage.meta.title" row (row 46), cloning the
#    layout/formatting of the last existing data row (row 45) and then
#    overwriting the three populated cells' values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: row 3 (homePage.meta.description) ---------------------------
$ws.Cells.Item(3, 2).Value = "轻盈好用的在线工具，无需下载即可免费使用，解决生活学习工作中的大小问题"
$ws.Cells.Item(3, 3).Value = "Effortless solutions at your fingertips. Streamline your workflow with our intuitive online toolbox. No installation. No hassle. Just pure productivity."

# C3 was previously blank (plain style); give it the same data-cell style
# used elsewhere in the column (e.g. C4) instead of the blank style.
$ws.Cells.Item(4, 3).Copy()
$ws.Cells.Item(3, 3).PasteSpecial(-4122)

# --- 3: new row 46 (aboutPage.meta.title) --------------------------------
$ws.Range("A45:E45").Copy($ws.Range("A46:E46"))
$ws.Rows.Item(46).RowHeight = $ws.Rows.Item(45).RowHeight

$ws.Cells.Item(46, 1).Value = "aboutPage.meta.title"
$ws.Cells.Item(46, 2).Value = "关于"
$ws.Cells.Item(46, 3).Value = "About"
